{"js": "// 1) Remove the bullet item \"Ing\u00e5 f\u00f6rlikning \u00e5 mina v\u00e4gnar\" from the\n//    \"R\u00e4tteg\u00e5ngsfullmakt\" scope list.\nconst bulletHits = context.document.body.search(\"Ing\u00e5 f\u00f6rlikning \u00e5 mina v\u00e4gnar\", { matchCase: true });\nbulletHits.load(\"text\");\nawait context.sync();\nif (bulletHits.items.length > 0) {\n  const bulletPara = bulletHits.items[0].paragraphs.getFirst();\n  bulletPara.delete();\n  await context.sync();\n}\n\n// 2) Remove special condition \"3. Fullmaktstagaren f\u00e5r ing\u00e5 f\u00f6rlikning \u00e5 mina\n//    v\u00e4gnar.\" entirely, and renumber the following condition from \"4.\" to \"3.\".\nconst cond3Hits = context.document.body.search(\n  \"3. Fullmaktstagaren f\u00e5r ing\u00e5 f\u00f6rlikning \u00e5 mina v\u00e4gnar.\",\n  { matchCase: true }\n);\ncond3Hits.load(\"text\");\nawait context.sync();\nif (cond3Hits.items.length > 0) {\n  const cond3Para = cond3Hits.items[0].paragraphs.getFirst();\n  cond3Para.delete();\n  await context.sync();\n}\n\nconst cond4Hits = context.document.body.search(\n  \"4. Jag f\u00f6rbinder mig att godk\u00e4nna de \u00e5tg\u00e4rder som fullmaktstagaren vidtar inom ramen f\u00f6r denna fullmakt.\",\n  { matchCase: true }\n);\ncond4Hits.load(\"text\");\nawait context.sync();\nif (cond4Hits.items.length > 0) {\n  cond4Hits.items[0].insertText(\n    \"3. Jag f\u00f6rbinder mig att godk\u00e4nna de \u00e5tg\u00e4rder som fullmaktstagaren vidtar inom ramen f\u00f6r denna fullmakt.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 3) Bump both \"Ort och datum: Eskilstuna den ______________ 2025\" lines to 2026.\nconst dateHits = context.document.body.search(\n  \"Ort och datum: Eskilstuna den ______________ 2025\",\n  { matchCase: true }\n);\ndateHits.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < dateHits.items.length; i++) {\n  dateHits.items[i].insertText(\n    \"Ort och datum: Eskilstuna den ______________ 2026\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 4) Drop one of the two blank paragraphs right after\n//    \"Joumana Alnablsi (19880102-5084)\" (collapse double blank to single blank).\nconst signatureHits = context.document.body.search(\n  \"Joumana Alnablsi (19880102-5084)\",\n  { matchCase: true }\n);\nsignatureHits.load(\"text\");\nawait context.sync();\nif (signatureHits.items.length > 0) {\n  const signaturePara = signatureHits.items[0].paragraphs.getFirst();\n  const blankAfterSignature = signaturePara.getNext();\n  blankAfterSignature.load(\"text\");\n  await context.sync();\n  if (blankAfterSignature.text === \"\") {\n    blankAfterSignature.delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the bullet item \"Ing\u00e5 f\u00f6rlikning \u00e5 mina v\u00e4gnar\" from the\n#    \"R\u00e4tteg\u00e5ngsfullmakt\" scope list.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Ing\u00e5 f\u00f6rlikning \u00e5 mina v\u00e4gnar\")\nif ($found) {\n    $bulletPara = $rng.Paragraphs(1)\n    $bulletPara.Range.Delete()\n}\n\n# 2) Remove special condition \"3. Fullmaktstagaren f\u00e5r ing\u00e5 f\u00f6rlikning \u00e5 mina\n#    v\u00e4gnar.\" entirely, and renumber the following condition from \"4.\" to \"3.\".\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"3. Fullmaktstagaren f\u00e5r ing\u00e5 f\u00f6rlikning \u00e5 mina v\u00e4gnar.\")\nif ($found2) {\n    $cond3Para = $rng2.Paragraphs(1)\n    $cond3Para.Range.Delete()\n}\n\n$rng3 = $d.Content\n$found3 = $rng3.Find.Execute(\"4. Jag f\u00f6rbinder mig att godk\u00e4nna de \u00e5tg\u00e4rder som fullmaktstagaren vidtar inom ramen f\u00f6r denna fullmakt.\")\nif ($found3) {\n    $cond4Para = $rng3.Paragraphs(1)\n    $cond4Para.Range.Text = \"3. Jag f\u00f6rbinder mig att godk\u00e4nna de \u00e5tg\u00e4rder som fullmaktstagaren vidtar inom ramen f\u00f6r denna fullmakt.\"\n}\n\n# 3) Bump both \"Ort och datum: Eskilstuna den ______________ 2025\" lines to 2026.\n$rng4 = $d.Content\n$rng4.Find.Execute(\"Eskilstuna den ______________ 2025\", $false, $false, $false, $false, $false, $true, \"wdFindContinue\", $false, \"Eskilstuna den ______________ 2026\", \"wdReplaceAll\")\n\n# 4) Drop one of the two blank paragraphs right after\n#    \"Joumana Alnablsi (19880102-5084)\" (collapse double blank to single blank).\n$rng5 = $d.Content\n$found5 = $rng5.Find.Execute(\"Joumana Alnablsi (19880102-5084)\")\nif ($found5) {\n    $signaturePara = $rng5.Paragraphs(1)\n    $blankAfterSignature = $signaturePara.Next()\n    if ($blankAfterSignature.Range.Text.Trim() -eq \"\") {\n        $blankAfterSignature.Range.Delete()\n    }\n}\n"}
